# Updated symbol list (crypto price/volume refresh) matching the commit:
# "Updated symbol list ... with GitHub Actions"
#
# Price (col D) and Volume(1h) (col E) are stored as plain text in this
# sheet, not as real numbers/percentages, so each numeric-looking value is
# written with a leading apostrophe to force text entry (same as typing
# '329.16 into a cell in Excel) instead of letting it auto-convert to a
# Number/Percentage. Coin name / link columns (B/C) are plain text already.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'329.16"
$ws.Range("E2").Value = "'1.84%"
$ws.Range("D3").Value = "'41.09"
$ws.Range("E3").Value = "'3.49%"
$ws.Range("D4").Value = "'5.620"
$ws.Range("E4").Value = "'-4.51%"
$ws.Range("D5").Value = "'0.08165"
$ws.Range("E5").Value = "'1.69%"
$ws.Range("D6").Value = "'2.037"
$ws.Range("E6").Value = "'5.69%"
$ws.Range("D7").Value = "'8.734"
$ws.Range("E7").Value = "'0.86%"
$ws.Range("D8").Value = "'4.534"
$ws.Range("E8").Value = "'-1.05%"
$ws.Range("D9").Value = "'2.949"
$ws.Range("E9").Value = "'0.15%"
$ws.Range("E10").Value = "'-1.57%"
$ws.Range("D11").Value = "'0.1256"
$ws.Range("E11").Value = "'-1.23%"
$ws.Range("D12").Value = "'0.1948"
$ws.Range("E12").Value = "'-1.06%"
$ws.Range("D13").Value = "'0.09287"
$ws.Range("E13").Value = "'0.29%"
$ws.Range("D14").Value = "'0.03739"
$ws.Range("E14").Value = "'5.81%"
$ws.Range("D16").Value = "'0.001302"
$ws.Range("E16").Value = "'0.73%"
$ws.Range("D17").Value = "'0.006174"
$ws.Range("E17").Value = "'1.57%"
$ws.Range("D18").Value = "'3.437"
$ws.Range("E18").Value = "'2.63%"
$ws.Range("E19").Value = "'-2.14%"
$ws.Range("D20").Value = "'8.263"
$ws.Range("E20").Value = "'-5.33%"
$ws.Range("D21").Value = "'0.1394"
$ws.Range("E21").Value = "'-1.74%"
$ws.Range("D22").Value = "'0.2653"
$ws.Range("E22").Value = "'10.18%"
$ws.Range("D23").Value = "'0.04419"
$ws.Range("E23").Value = "'0.22%"
$ws.Range("D24").Value = "'0.001272"
$ws.Range("E24").Value = "'0.92%"
$ws.Range("D25").Value = "'0.004288"
$ws.Range("E25").Value = "'-1.95%"
$ws.Range("E26").Value = "'3.75%"
$ws.Range("D39").Value = "'0.02771"
$ws.Range("E39").Value = "'13.38%"
$ws.Range("D40").Value = "'0.05413"
$ws.Range("E40").Value = "'3.19%"
$ws.Range("D41").Value = "'0.007670"
$ws.Range("E41").Value = "'2.30%"
$ws.Range("D42").Value = "'0.009408"
$ws.Range("E42").Value = "'-0.64%"
$ws.Range("D43").Value = "'0.1414"
$ws.Range("E43").Value = "'0.66%"
$ws.Range("D44").Value = "'0.002134"
$ws.Range("E44").Value = "'0.70%"
$ws.Range("D45").Value = "'0.01142"
$ws.Range("E45").Value = "'15.81%"
$ws.Range("D46").Value = "'0.00006895"
$ws.Range("E46").Value = "'2.34%"
$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("E47").Value = "'0.23%"
$ws.Range("B48").Value = "CoinbaseStockToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D48").Value = "'0.002285"
$ws.Range("E48").Value = "'60.59%"
$ws.Range("B49").Value = "BOLO"
$ws.Range("C49").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D49").Value = "'0.003552"
$ws.Range("E49").Value = "'18.40%"
$ws.Range("E50").Value = "'0.23%"
$ws.Range("E51").Value = "'0.23%"
